$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "28-jul" column (AN) header, mirroring the existing "27-jul" (AM) column
$ws.Range("AN1").Value = "28-jul"

# New daily values for each product line row
$ws.Range("AN2").Value = 0
$ws.Range("AN3").Value = 18.059946848652984
$ws.Range("AN4").Value = 15.628244088871613
$ws.Range("AN5").Value = 22.058646804729133
$ws.Range("AN6").Value = 0
$ws.Range("AN7").Value = 4.5434088342255592
$ws.Range("AN8").Value = 5.8319774095429713
$ws.Range("AN9").Value = 16.812588227818729
$ws.Range("AN10").Value = 16.315969771958027
$ws.Range("AN11").Value = 8.8123318708146101
$ws.Range("AN12").Value = 0
$ws.Range("AN13").Value = 12.971410086591584
$ws.Range("AN14").Value = 0
$ws.Range("AN15").Value = 0
$ws.Range("AN16").Value = 2.4721278280188441
$ws.Range("AN17").Value = 0
$ws.Range("AN18").Value = 0

# Move the active selection, as recorded after the edit
$ws.Range("AO4").Select()
